$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 654, pushing existing rows 654:695 down to 655:696
$ws.Rows.Item(654).Insert()

# Populate the newly inserted row 654 with the new record.
# Column A holds a date-formatted string ("2026/01/17"); force it to be
# stored as plain text (matching the rest of the sheet) instead of letting
# Excel auto-convert it to a date serial number.
$dateCell = $ws.Cells.Item(654, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/01/17"
$dateCell.ClearFormats()

$ws.Cells.Item(654, 2).Value = "土"
$ws.Cells.Item(654, 3).Value = 5
$ws.Cells.Item(654, 4).Value = 201
